$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row 1 labels: _old -> _FV2304, _new -> _FV2310
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cur = $cell.Value()
    $cell.Value = ($cur -replace "_old$", "_FV2304")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cur = $cell.Value()
    $cell.Value = ($cur -replace "_new$", "_FV2310")
}
